# "New tenant support in live" - append the 176* regression-run rows that
# landed for the new tenant onto the AMSIN, BETA and AMS sprint-history
# sheets, and backfill the formatting on AMSIN!A58:G58 (the previous last
# row) so it matches the rest of the table.

$wb = $excel.ActiveWorkbook

function Write-HistoryRow {
    param($ws, $row, $runDate, $runTime, $sprintName, $totalCases, $passCases, $failCases, $timeTaken, $dateStyleSourceRow)

    # Column A ("Run Date") is stored as literal text ("2023-04-18", ...),
    # not a real date - format the cell as Text first so assigning a
    # date-shaped string doesn't get auto-converted into a date serial.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $runDate
    $ws.Cells.Item($row, 1).NumberFormat = "general"

    # Column B ("Run Time") keeps the existing date/time display format -
    # clone it off a neighbouring cell that already has it, then overwrite
    # the value.
    $ws.Range("B" + $dateStyleSourceRow).Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4122)
    $ws.Cells.Item($row, 2).Value = $runTime

    $ws.Cells.Item($row, 3).Value = $sprintName
    $ws.Cells.Item($row, 4).Value = $totalCases
    $ws.Cells.Item($row, 5).Value = $passCases
    $ws.Cells.Item($row, 6).Value = $failCases
    $ws.Cells.Item($row, 7).Value = $timeTaken
}

# --- AMSIN: fix up row 58's formatting, then append rows 59-60 ---
$wsAmsin = $wb.Worksheets.Item("AMSIN")

# Row 58 itself picks up the same (blank/general) styling as the rest of
# the table and its run-time gets re-stamped to millisecond precision.
$wsAmsin.Range("A57:G57").Copy()
$wsAmsin.Range("A58:G58").PasteSpecial(-4122)
$wsAmsin.Cells.Item(58, 2).Value = 45034.61307858796

Write-HistoryRow $wsAmsin 59 "2023-04-19" 45035.71441813657 "176scndcyc" 119 119 0 3.06 57
Write-HistoryRow $wsAmsin 60 "2023-04-20" 45036.42718866898 "176fnlruntest" 119 119 0 3.25 57

# --- BETA: append row 30 ---
$wsBeta = $wb.Worksheets.Item("BETA")

Write-HistoryRow $wsBeta 30 "2023-04-20" 45036.53269996528 "176beta" 119 119 0 2.82 29

# --- AMS: append row 31 ---
$wsAms = $wb.Worksheets.Item("AMS")

Write-HistoryRow $wsAms 31 "2023-05-08" 45054.55315972622 "176htfxtrl" 119 119 0 2.61 30
